$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price column (D) - force text to avoid numeric auto-conversion
$dValues = [ordered]@{
    "D2" = "69.485.42"
    "D3" = "3.771.69"
    "D5" = "615.71"
    "D6" = "178.08"
    "D7" = "3.771.91"
    "D9" = "0.528"
    "D11" = "6.68"
    "D13" = "40.06"
    "D15" = "4.396.24"
    "D16" = "3.770.45"
    "D17" = "69.517.68"
    "D20" = "509.68"
    "D21" = "16.34"
    "D22" = "9.36"
    "D24" = "2.52"
    "D25" = "86.41"
    "D28" = "10.60"
    "D32" = "7.99"
    "D33" = "30.68"
    "D35" = "1.00"
    "D37" = "6.14"
    "D39" = "0.340"
    "D40" = "453.08"
    "D42" = "49.87"
    "D43" = "2.99"
    "D44" = "44.68"
    "D45" = "8.56"
    "D46" = "2.958.05"
    "D47" = "0.0359"
    "D49" = "139.08"
    "D50" = "27.22"
    "D51" = "2.46"
}
foreach ($ref in $dValues.Keys) {
    $ws.Range($ref).Formula = "'" + $dValues[$ref]
    $ws.Range($ref).Style = "Normal"
}

# Update Volume(1h) column (E)
$eValues = [ordered]@{
    "E2" = "  -0.72%  "
    "E3" = "  -0.37%  "
    "E4" = "  +0.02%  "
    "E5" = "  -1.01%  "
    "E6" = "  +0.31%  "
    "E7" = "  -0.27%  "
    "E8" = "  +0.00%  "
    "E9" = "  -1.46%  "
    "E10" = "  -2.53%  "
    "E11" = "  +5.89%  "
    "E12" = "  -1.71%  "
    "E13" = "  -2.29%  "
    "E14" = "  -3.73%  "
    "E15" = "  -0.35%  "
    "E16" = "  -0.30%  "
    "E17" = "  -0.72%  "
    "E18" = "  -1.20%  "
    "E19" = "  -3.50%  "
    "E20" = "  -0.24%  "
    "E21" = "  -2.98%  "
    "E22" = "  -1.64%  "
    "E23" = "  +0.03%  "
    "E24" = "  +1.05%  "
    "E25" = "  -1.60%  "
    "E26" = "  -2.46%  "
    "E27" = "  -2.48%  "
    "E28" = "  -3.59%  "
    "E29" = "  +0.16%  "
    "E30" = "  +0.47%  "
    "E31" = "  +2.97%  "
    "E32" = "  +2.94%  "
    "E33" = "  -2.24%  "
    "E34" = "  -1.23%  "
    "E35" = "  -0.02%  "
    "E36" = "  -0.77%  "
    "E37" = "  -1.08%  "
    "E38" = "  +3.60%  "
    "E39" = "  +2.19%  "
    "E40" = "  +8.45%  "
    "E41" = "  -2.95%  "
    "E42" = "  -2.19%  "
    "E43" = "  +5.88%  "
    "E44" = "  -0.62%  "
    "E45" = "  -2.19%  "
    "E46" = "  -2.49%  "
    "E47" = "  -1.26%  "
    "E48" = "  +0.03%  "
    "E49" = "  +0.11%  "
    "E50" = "  -0.71%  "
    "E51" = "  -1.16%  "
}
foreach ($ref in $eValues.Keys) {
    $ws.Range($ref).Value = $eValues[$ref]
}

Write-Host "Applied $($dValues.Count) price updates and $($eValues.Count) volume updates"